$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Add narrative "It was an amazing class!" to the feedbacks row
#    (feedback_id 541, rating 5) at row 71.
$ws.Cells.Item(71, 3).Value = "It was an amazing class!"

# 2. Add description "Payment for Course" to existing issues row (row 34)
$ws.Cells.Item(34, 5).Value = "Payment for Course"

# 3. Insert a new row at row 35 (pushes everything below down by one),
#    then populate it as a new "issues" record.
$ws.Rows(35).Insert()
$ws.Cells.Item(35, 1).Value = 4322
$ws.Cells.Item(35, 2).Value = 203
$ws.Cells.Item(35, 4).Value = 8521
$ws.Cells.Item(35, 5).Value = "Payment for Course"
$ws.Cells.Item(35, 6).Value = 201

# 4. Update the active selection to match (G34).
$ws.Activate()
$ws.Range("G34").Select()
